$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 88) from DGS's 2021/10/01 report.
$row = 88

# Column A holds the date label as text (matching the existing rows,
# which are shared strings, not real dates). Entering the literal text
# into a cell formatted as a date would normally get auto-converted to
# a date serial number, so build it via a text formula and then convert
# the formula result back to a plain value in place.
$ws.Cells.Item($row, 1).Formula = '="2021/10/01"'
$ws.Range("A" + $row).Copy()
$ws.Range("A" + $row).PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 2).Value = 101.7
$ws.Cells.Item($row, 3).Value = 103
$ws.Cells.Item($row, 4).Value = 0.89
$ws.Cells.Item($row, 5).Value = 0.89

$ws.Cells.Item($row + 1, 1).Select()
